# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets.
# F2: 8658 -> 8663
# F4: 402  -> 403

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 8663
    $ws.Range("F4").Value = 403
}
